# Updated cryptos list with latest data
# Applies per-cell updates for Price (column D) and Volume(1h) (column E),
# plus two rows that had their coin identity swapped (rows 26/27) and
# one row whose coin was fully replaced (row 51: Aptos -> SynthetixNetwork).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to stay a text value even when the string looks like
    # a number (e.g. "1.000", "29.357.13"), matching the original inline
    # string cells, then restore the default General display format.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
}


Set-TextValue $ws.Range("D2") '29.357.13'
Set-TextValue $ws.Range("E2") '  +0.35%  '
Set-TextValue $ws.Range("D3") '1.872.06'
Set-TextValue $ws.Range("E3") '  +0.67%  '
Set-TextValue $ws.Range("E4") '  -0.01%  '
Set-TextValue $ws.Range("D5") '0.7082'
Set-TextValue $ws.Range("E5") '  +0.80%  '
Set-TextValue $ws.Range("D6") '238.73'
Set-TextValue $ws.Range("E6") '  +0.45%  '
Set-TextValue $ws.Range("D7") '1.000'
Set-TextValue $ws.Range("E7") '  -0.01%  '
Set-TextValue $ws.Range("D8") '0.07784'
Set-TextValue $ws.Range("E8") '  -5.56%  '
Set-TextValue $ws.Range("D9") '0.3064'
Set-TextValue $ws.Range("E9") '  +0.93%  '
Set-TextValue $ws.Range("D10") '25.18'
Set-TextValue $ws.Range("E10") '  +8.26%  '
Set-TextValue $ws.Range("D11") '0.08195'
Set-TextValue $ws.Range("E11") '  +0.18%  '
Set-TextValue $ws.Range("D12") '1.883.73'
Set-TextValue $ws.Range("E12") '  +0.91%  '
Set-TextValue $ws.Range("D13") '5.242'
Set-TextValue $ws.Range("E13") '  +1.23%  '
Set-TextValue $ws.Range("D14") '0.7206'
Set-TextValue $ws.Range("E14") '  +0.76%  '
Set-TextValue $ws.Range("D15") '89.21'
Set-TextValue $ws.Range("E15") '  +0.06%  '
Set-TextValue $ws.Range("D16") '29.398.92'
Set-TextValue $ws.Range("E16") '  +0.41%  '
Set-TextValue $ws.Range("D17") '5.810'
Set-TextValue $ws.Range("E17") '  +0.51%  '
Set-TextValue $ws.Range("D18") '242.24'
Set-TextValue $ws.Range("E18") '  +2.21%  '
Set-TextValue $ws.Range("D19") '0.000007832'
Set-TextValue $ws.Range("E19") '  -0.13%  '
Set-TextValue $ws.Range("D20") '13.27'
Set-TextValue $ws.Range("E20") '  -0.75%  '
Set-TextValue $ws.Range("D21") '2.125.38'
Set-TextValue $ws.Range("E21") '  -0.39%  '
Set-TextValue $ws.Range("D22") '1.000'
Set-TextValue $ws.Range("E22") '  +0.01%  '
Set-TextValue $ws.Range("D23") '1.001'
Set-TextValue $ws.Range("E23") '  -0.04%  '
Set-TextValue $ws.Range("D24") '7.714'
Set-TextValue $ws.Range("E24") '  +3.43%  '
Set-TextValue $ws.Range("D25") '162.26'
Set-TextValue $ws.Range("B26") 'Stellar'
Set-TextValue $ws.Range("C26") 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D26") '0.1461'
Set-TextValue $ws.Range("E26") '  +1.31%  '
Set-TextValue $ws.Range("B27") 'Cosmos'
Set-TextValue $ws.Range("C27") 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D27") '8.960'
Set-TextValue $ws.Range("E27") '  -0.24%  '
Set-TextValue $ws.Range("D28") '18.16'
Set-TextValue $ws.Range("E28") '  +0.31%  '
Set-TextValue $ws.Range("D29") '1.925'
Set-TextValue $ws.Range("E29") '  -2.14%  '
Set-TextValue $ws.Range("D30") '1.367'
Set-TextValue $ws.Range("E30") '  -5.07%  '
Set-TextValue $ws.Range("D31") '1.484'
Set-TextValue $ws.Range("E31") '  +0.18%  '
Set-TextValue $ws.Range("D32") '4.311'
Set-TextValue $ws.Range("E32") '  -2.37%  '
Set-TextValue $ws.Range("D33") '4.060'
Set-TextValue $ws.Range("E33") '  +0.06%  '
Set-TextValue $ws.Range("D34") '0.05213'
Set-TextValue $ws.Range("E34") '  +0.02%  '
Set-TextValue $ws.Range("E35") '  +1.84%  '
Set-TextValue $ws.Range("E36") '  +1.56%  '
Set-TextValue $ws.Range("E37") '  +0.47%  '
Set-TextValue $ws.Range("E38") '  +0.31%  '
Set-TextValue $ws.Range("E39") '  +0.28%  '
Set-TextValue $ws.Range("D40") '2.701'
Set-TextValue $ws.Range("E40") '  -0.91%  '
Set-TextValue $ws.Range("D41") '1.174.46'
Set-TextValue $ws.Range("E41") '  +3.65%  '
Set-TextValue $ws.Range("D42") '0.9170'
Set-TextValue $ws.Range("D43") '5.996'
Set-TextValue $ws.Range("E43") '  +0.66%  '
Set-TextValue $ws.Range("D44") '0.4292'
Set-TextValue $ws.Range("D45") '71.42'
Set-TextValue $ws.Range("E45") '  +0.98%  '
Set-TextValue $ws.Range("D46") '1.000'
Set-TextValue $ws.Range("E46") '  +0.05%  '
Set-TextValue $ws.Range("D47") '102.42'
Set-TextValue $ws.Range("E47") '  -0.20%  '
Set-TextValue $ws.Range("D48") '0.5299'
Set-TextValue $ws.Range("E48") '  -2.22%  '
Set-TextValue $ws.Range("D49") '1.759'
Set-TextValue $ws.Range("E49") '  -0.78%  '
Set-TextValue $ws.Range("D50") '9.200'
Set-TextValue $ws.Range("E50") '  +0.20%  '
Set-TextValue $ws.Range("B51") 'SynthetixNetwork'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
Set-TextValue $ws.Range("D51") '2.860'
Set-TextValue $ws.Range("E51") '  +3.21%  '
